$wb = $excel.ActiveWorkbook

# Add the new worksheet "tc002" right after the existing "tc001" sheet
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "tc002"

# Populate data
$ws2.Range("A1").Value = "projectName"
$ws2.Range("A2").Value = "SET- DRV"

# Column width (engine rounds to nearest 1/6-char increment; 22.8 -> stored width 23.667,
# the closest achievable value to the authored 23.625)
$ws2.Columns.Item(1).ColumnWidth = 22.8

# Select A7 and make tc002 the active/visible sheet
$ws2.Range("A7").Select()
$ws2.Activate()
